# Update filtered_output.xlsx: refresh the Neg_Change and Pos_Change
# sheets with newly filtered market data (rows differ from the old data,
# and row counts shrink from 14/15 data rows down to 11/12 data rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Neg_Change
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Neg_Change")

$negData = @(
    @("RELIANCE",   1510,   1520,    1498,    1505.1,  11199340, 27417125, -0.5915202633390627),
    @("ITC",        343,    344.55,  339.3,   341.75,  29440734, 72582738, -0.5943838051411067),
    @("HDFCLIFE",   778,    779.9,   769.5,   773.7,   1349921,  2900960,  -0.5346640422480834),
    @("BHARTIARTL", 2095.6, 2103.8,  2077,    2084,    4171371,  9119265,  -0.5425759641813238),
    @("SIEMENS",    3120,   3138.4,  3090,    3131.8,  190568,   374745,   -0.4914728682170543),
    @("JINDALSTEL", 1080,   1087,    1067.6,  1074,    553494,   1100849,  -0.4972116975170982),
    @("MOTHERSON",  120.7,  121.27,  118.73,  119.17,  9250970,  19478525, -0.525068248237482),
    @("BANKINDIA",  150.66, 151.97,  149.71,  151.35,  5035903,  11879226, -0.5760748217097645),
    @("BDL",        1531,   1556.5,  1526.5,  1540,    2046391,  4768472,  -0.5708497397069753),
    @("COLPAL",     2089,   2103.5,  2061.1,  2075,    199951,   459339,   -0.5646984035755728),
    @("SBICARD",    901.05, 901.45,  881.35,  883.45,  829925,   1737746,  -0.5224129418223377)
)

$destRows1 = $negData.Count + 1

# delete rows that won't be needed anymore (old sheet had data through row 14)
if ($ws1.UsedRange.Rows.Count -gt $destRows1) {
    $ws1.Range("A" + ($destRows1 + 1) + ":I14").EntireRow.Delete()
}

$r = 2
foreach ($row in $negData) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $ws1.Cells.Item($r, 9).Value = $row[0]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: Pos_Change
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Pos_Change")

$posData = @(
    @("COALINDIA",  426.95, 432.45,  423.65,  431.4,   8920605,   6308113,   0.4141479393282904),
    @("BRITANNIA",  6200,   6271,    6157.5,  6175,    437361,    300112,    0.4573259316521832),
    @("PIDILITIND", 1510,   1517.5,  1496.6,  1511,    623936,    393706,    0.5847764575597019),
    @("INDHOTEL",   720,    721.05,  703.3,   714,     3998970,   2525206,   0.5836212966387693),
    @("YESBANK",    22.89,  23.57,   22.72,   23.46,   188593277, 121289455, 0.5549025016230801),
    @("DIXON",      11640,  11970,   11480,   11780,   837880,    541945,    0.5460609471440829),
    @("UNOMINDA",   1321,   1334.9,  1302.3,  1316.2,  607276,    396153,    0.5329329829636529),
    @("SUZLON",     53.3,   53.69,   52.67,   52.96,   48744390,  31119731,  0.5663499790534822),
    @("PETRONET",   293.5,  299.8,   291.1,   292.4,   6494904,   4129368,   0.5728566695920538),
    @("IGL",        191,    191.42,  186.65,  187.2,   1304820,   834609,    0.5633907614224146),
    @("JUBLFOOD",   547.8,  556.15,  535.05,  536,     3508452,   2208300,   0.5887569623692434),
    @("CYIENT",     1132.2, 1184.5,  1132.2,  1171.2,  379839,    263020,    0.4441449319443388)
)

$destRows2 = $posData.Count + 1

# delete rows that won't be needed anymore (old sheet had data through row 15)
if ($ws2.UsedRange.Rows.Count -gt $destRows2) {
    $ws2.Range("A" + ($destRows2 + 1) + ":I15").EntireRow.Delete()
}

$r = 2
foreach ($row in $posData) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $ws2.Cells.Item($r, 8).Value = $row[7]
    $ws2.Cells.Item($r, 9).Value = $row[0]
    $r = $r + 1
}
